$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.150.80"
$ws.Range("E2").Value = "  +0.47%  "

# Row 3
$ws.Range("D3").Value = "2.546.87"
$ws.Range("E3").Value = "  -2.30%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "590.08"
$ws.Range("E5").Value = "  +0.79%  "

# Row 6
$ws.Range("D6").Value = "173.32"
$ws.Range("E6").Value = "  +5.07%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  +0.46%  "

# Row 9
$ws.Range("D9").Value = "2.546.09"
$ws.Range("E9").Value = "  -2.25%  "

# Row 10
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  -0.88%  "

# Row 11
$ws.Range("E11").Value = "  +1.85%  "

# Row 12
$ws.Range("D12").Value = "5.14"
$ws.Range("E12").Value = "  -0.66%  "

# Row 13
$ws.Range("D13").Value = "0.347"
$ws.Range("E13").Value = "  -4.72%  "

# Row 14
$ws.Range("D14").Value = "26.88"
$ws.Range("E14").Value = "  -0.75%  "

# Row 15
$ws.Range("D15").Value = "3.011.60"
$ws.Range("E15").Value = "  -2.38%  "

# Row 16
$ws.Range("E16").Value = "  -0.88%  "

# Row 17
$ws.Range("D17").Value = "66.963.84"
$ws.Range("E17").Value = "  +0.31%  "

# Row 18
$ws.Range("D18").Value = "2.552.75"
$ws.Range("E18").Value = "  -1.99%  "

# Row 19
$ws.Range("D19").Value = "8.02"
$ws.Range("E19").Value = "  +3.34%  "

# Row 20
$ws.Range("D20").Value = "11.30"
$ws.Range("E20").Value = "  -2.74%  "

# Row 21
$ws.Range("D21").Value = "355.17"
$ws.Range("E21").Value = "  +0.47%  "

# Row 22
$ws.Range("D22").Value = "4.19"
$ws.Range("E22").Value = "  -1.28%  "

# Row 23
$ws.Range("D23").Value = "4.64"
$ws.Range("E23").Value = "  +0.88%  "

# Row 24
$ws.Range("E24").Value = "  +6.48%  "

# Row 25
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("D26").Value = "69.78"
$ws.Range("E26").Value = "  +0.86%  "

# Row 27
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  -3.99%  "

# Row 28
$ws.Range("E28").Value = "  -2.48%  "

# Row 29
$ws.Range("D29").Value = "0.994"
$ws.Range("E29").Value = "  -0.43%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0986"
$ws.Range("E30").Value = "  +0.18%  "

# Row 31
$ws.Range("D31").Value = "532.79"
$ws.Range("E31").Value = "  -0.93%  "

# Row 32
$ws.Range("D32").Value = "8.17"
$ws.Range("E32").Value = "  +0.89%  "

# Row 33
$ws.Range("E33").Value = "  +1.45%  "

# Row 34
$ws.Range("E34").Value = "  -0.46%  "

# Row 36
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("E37").Value = "  +0.33%  "

# Row 38
$ws.Range("D38").Value = "157.23"
$ws.Range("E38").Value = "  -0.47%  "

# Row 39
$ws.Range("D39").Value = "18.68"

# Row 40
$ws.Range("E40").Value = "  +1.19%  "

# Row 41
$ws.Range("D41").Value = "0.356"
$ws.Range("E41").Value = "  -1.48%  "

# Row 42
$ws.Range("E42").Value = "  +0.24%  "

# Row 43
$ws.Range("D43").Value = "5.14"
$ws.Range("E43").Value = "  +1.03%  "

# Row 44
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.51"
$ws.Range("E45").Value = "  +5.01%  "

# Row 46
$ws.Range("E46").Value = "  -1.19%  "

# Row 47
$ws.Range("D47").Value = "149.77"
$ws.Range("E47").Value = "  -0.22%  "

# Row 48
$ws.Range("D48").Value = "0.561"
$ws.Range("E48").Value = "  -1.87%  "

# Row 49
$ws.Range("E49").Value = "  -4.60%  "

# Row 50
$ws.Range("D50").Value = "3.69"
$ws.Range("E50").Value = "  -1.05%  "

# Row 51
$ws.Range("D51").Value = "1.70"
$ws.Range("E51").Value = "  +0.54%  "
